# "Update countries & provincias Spain"
#
# Refreshes the COVID-19 country table on sheet "Pais":
#  - bumps the "last updated" timestamp (A1)
#  - updates Casos totales / Nuevos casos / Casos activos / Recuperados /
#    Casos criticos / Muertes hoy / Muertes (cols B-H) for the countries
#    whose figures changed in this refresh
#  - two pairs of neighbouring countries swapped rank because one of them
#    overtook the other in "Casos totales": Emiratos Arabes Unidos now
#    outranks Polonia (rows 43/44), and Dinamarca now outranks Hungria
#    (rows 78/79). The country names in column A are updated accordingly
#    while the newly-reported figures follow the country that moved up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 2 de Octubre de 2020 a las 14:18"

# Row 40 - Kuwait
$ws.Range("B40").Value = 106087
$ws.Range("C40").Value = 411
$ws.Range("D40").Value = 97898
$ws.Range("E40").Value = 7574
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 615

# Row 43 - was Polonia, now Emiratos Arabes Unidos (moved up, new figures)
$ws.Range("A43").Value = "Emiratos Arabes Unidos"
$ws.Range("B43").Value = 96529
$ws.Range("C43").Value = 1181
$ws.Range("D43").Value = 86071
$ws.Range("E43").Value = 10034
$ws.Range("G43").Value = 3
$ws.Range("H43").Value = 424

# Row 44 - was Emiratos Arabes Unidos, now Polonia (moved down one spot, same figures as old row 43)
$ws.Range("A44").Value = "Polonia"
$ws.Range("B44").Value = 95773
$ws.Range("C44").Value = 2292
$ws.Range("D44").Value = 71353
$ws.Range("E44").Value = 21850
$ws.Range("G44").Value = 27
$ws.Range("H44").Value = 2570

# Row 45 - Suecia
$ws.Range("B45").Value = 94283
$ws.Range("G45").Value = 5
$ws.Range("H45").Value = 5895

# Row 49 - Nepal
$ws.Range("B49").Value = 82450
$ws.Range("C49").Value = 2722
$ws.Range("D49").Value = 60696
$ws.Range("E49").Value = 21234
$ws.Range("G49").Value = 11
$ws.Range("H49").Value = 520

# Row 61 - Suiza
$ws.Range("D61").Value = 45800
$ws.Range("E61").Value = 6509

# Row 70 - Estado de Palestina
$ws.Range("B70").Value = 40766
$ws.Range("C70").Value = 444
$ws.Range("D70").Value = 32944
$ws.Range("E70").Value = 7501
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = 321

# Row 71 - Azerbaiyan
$ws.Range("B71").Value = 40453
$ws.Range("C71").Value = 144
$ws.Range("D71").Value = 38217
$ws.Range("E71").Value = 1642
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 594

# Row 78 - was Hungria, now Dinamarca (moved up, new figures)
$ws.Range("A78").Value = "Dinamarca"
$ws.Range("B78").Value = 28932
$ws.Range("C78").Value = 536
$ws.Range("D78").Value = 21824
$ws.Range("E78").Value = 6456
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 652

# Row 79 - was Dinamarca, now Hungria (moved down one spot, same figures as old row 78)
$ws.Range("A79").Value = "Hungria"
$ws.Range("B79").Value = 28631
$ws.Range("C79").Value = 1322
$ws.Range("D79").Value = 6349
$ws.Range("E79").Value = 21484
$ws.Range("G79").Value = 17
$ws.Range("H79").Value = 798

# Row 88 - Republica de Macedonia
$ws.Range("B88").Value = 18363
$ws.Range("C88").Value = 225
$ws.Range("D88").Value = 15164
$ws.Range("E88").Value = 2450
$ws.Range("G88").Value = 6
$ws.Range("H88").Value = 749

# Row 102 - Consejo Danes para los Refugiados
$ws.Range("B102").Value = 10729
$ws.Range("C102").Value = 44
$ws.Range("D102").Value = 10183
$ws.Range("E102").Value = 274

# Row 112 - Uganda
$ws.Range("B112").Value = 8491
$ws.Range("C112").Value = 204
$ws.Range("D112").Value = 4470
$ws.Range("E112").Value = 3942
$ws.Range("G112").Value = 4
$ws.Range("H112").Value = 79

# Row 148 - Islandia
$ws.Range("B148").Value = 2809
$ws.Range("C148").Value = 40
$ws.Range("D148").Value = 2194
$ws.Range("E148").Value = 605

